# Generate Report for Handback
#
# Overview sheet: the per-language "Latest HO Xliff Generate Date" summary
#   (columns E/F) flips from "Ready for handoff" to
#   "Handed back: in sync with en-US" (same shared string is used by the
#   per-file "Status" column on the language sheets, so those pick up the
#   new text automatically once every referencing cell is rewritten).
# Per-language sheets (zh-cn / de-de): populate the "Latest Target File"
#   (I) and "Latest Handback File" (J) columns for each file row, stamp a
#   "Latest Handback DateTime" (K), and widen a few columns to fit the
#   newly-populated data.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$colWidthWide = 29.9777047293527 - (5/6)     # -> stored width 30 (closest grid point to 29.9777...)
$colWidthForty = 40 - (5/6)                  # -> stored width 40 exactly

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

foreach ($addr in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($addr)
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = $colWidthWide   # column E
$wsOverview.Columns.Item(6).ColumnWidth = $colWidthWide   # column F

# ---------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------
$langSheets = @(
    @{ Name = "zh-cn"; HandbackDate = "2016-09-04 05:08:29" },
    @{ Name = "de-de"; HandbackDate = "2016-09-04 05:08:37" }
)

foreach ($info in $langSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Collect the existing hyperlink addresses (keyed by anchor cell) before
    # adding new ones.
    $hyperlinkByRef = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $refAddr = $hl.Range.Address()
        $hyperlinkByRef[$refAddr] = $hl.Address
    }

    for ($row = 2; $row -le 3; $row++) {
        $statusCell = $ws.Cells.Item($row, 3)   # column C - Status
        if ($statusCell.Text -eq $oldStatus) {
            $statusCell.Value = $newStatus
        }

        $srcRef = "`$A`$$row"
        $srcUrl = $hyperlinkByRef[$srcRef]
        $srcName = $ws.Cells.Item($row, 1).Text       # column A - Source File Name (md)
        $handoffFile = $ws.Cells.Item($row, 7).Text   # column G - Latest Handoff File (xlf)

        # Column I - Latest Target File: hyperlink to the same source .md
        $targetCell = $ws.Cells.Item($row, 9)
        $ws.Hyperlinks.Add($targetCell, $srcUrl, "", "", $srcName) | Out-Null

        # Column J - Latest Handback File: same name as the handoff xlf
        $ws.Cells.Item($row, 10).Value = $handoffFile

        # Column K - Latest Handback DateTime
        $ws.Cells.Item($row, 11).Value = $info.HandbackDate
    }

    $ws.Columns.Item(3).ColumnWidth = $colWidthWide     # column C
    $ws.Columns.Item(9).ColumnWidth = $colWidthForty     # column I
    $ws.Columns.Item(10).ColumnWidth = $colWidthForty    # column J
}
